$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.023.72'
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").Value = '1.910.29'
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.16'
$ws.Range("E5").Value = '  +1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4802'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3812'
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07361'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9336'
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.81'
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '1.897.71'
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.497'
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.631'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.98'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").Value = '28.046.49'
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.78'
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.172'
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").Value = '2.127.58'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.52'
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.914'
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.49'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.136'
$ws.Range("E28").Value = '  +5.95%  '
$ws.Range("E29").Value = '  +1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.966'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08947'
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.300'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.265'
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7786'
$ws.Range("E34").Value = '  +3.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.681'
$ws.Range("E35").Value = '  +1.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.621'
$ws.Range("E36").Value = '  -4.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02055'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.113'
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05315'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5491'
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.019'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.483'
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.70'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4829'
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '107.95'
$ws.Range("E47").Value = '  +4.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.003'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.651'
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.93'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06073'
$ws.Range("E51").Value = '  -0.08%  '
